$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.413.13'
$ws.Range("E2").Value = '  +5.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.051.38'
$ws.Range("E3").Value = '  +3.49%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.66'
$ws.Range("E5").Value = '  +3.05%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.652'
$ws.Range("E6").Value = '  +2.60%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.46'
$ws.Range("E7").Value = '  +13.69%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +6.52%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.17'
$ws.Range("E10").Value = '  +1.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0767'
$ws.Range("E11").Value = '  +4.56%  '

# Row 12
$ws.Range("E12").Value = '  +1.51%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.922'
$ws.Range("E13").Value = '  -2.54%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.86'
$ws.Range("E14").Value = '  +3.04%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.59'
$ws.Range("E15").Value = '  +26.31%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.352.53'
$ws.Range("E16").Value = '  +3.55%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.58'
$ws.Range("E17").Value = '  +5.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.052.77'
$ws.Range("E18").Value = '  +3.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.328.11'
$ws.Range("E19").Value = '  +4.98%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.52'
$ws.Range("E20").Value = '  +3.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0874'
$ws.Range("E21").Value = '  +3.70%  '

# Row 22
$ws.Range("E22").Value = '  +6.22%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.83'
$ws.Range("E23").Value = '  +2.99%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.67'
$ws.Range("E24").Value = '  +4.66%  '

# Row 25
$ws.Range("E25").Value = '  -0.07%  '

# Row 26
$ws.Range("E26").Value = '  +4.87%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  +11.09%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.86'
$ws.Range("E28").Value = '  -1.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.97'
$ws.Range("E29").Value = '  +4.13%  '

# Row 30
$ws.Range("E30").Value = '  +26.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.26'
$ws.Range("E31").Value = '  +8.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.122'
$ws.Range("E32").Value = '  +2.85%  '

# Row 33
$ws.Range("E33").Value = '  +9.49%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0627'
$ws.Range("E34").Value = '  +5.83%  '

# Row 35
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.68'
$ws.Range("E35").Value = '  +8.43%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").Value = '  +0.73%  '

# Row 37
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("E38").Value = '  +3.87%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.05'
$ws.Range("E39").Value = '  +15.49%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.03'
$ws.Range("E40").Value = '  +35.04%  '

# Row 41
$ws.Range("E41").Value = '  +17.93%  '

# Row 42
$ws.Range("E42").Value = '  +4.48%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.24'
$ws.Range("E43").Value = '  +1.93%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.52'
$ws.Range("E44").Value = '  +9.40%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.15'
$ws.Range("E45").Value = '  +5.91%  '

# Row 46
$ws.Range("E46").Value = '  +3.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '96.44'
$ws.Range("E47").Value = '  +5.14%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.92'
$ws.Range("E48").Value = '  +3.16%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.421.74'
$ws.Range("E49").Value = '  +4.03%  '

# Row 50
$ws.Range("E50").Value = '  +1.96%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.58'
$ws.Range("E51").Value = '  -1.18%  '

